$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.1145763333333333
$ws.Range("H2").Value = 0.343729
$ws.Range("I2").Value = 0.006557053879060051
$ws.Range("J2").Value = 0.006557053879060051
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 34.07074633333333
$ws.Range("N2").Value = 102.212239
$ws.Range("O2").Value = 0.5171464495142372
$ws.Range("P2").Value = 0.5171464495142373
$ws.Range("Q2").Value = 3.903701188803443
$ws.Range("R2").Value = 35.13331069923099
$ws.Range("S2").Value = 0.003390957132829461
$ws.Range("T2").Value = 0.003390957132829462

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.1145763333333333
$ws.Range("H3").Value = 0.343729
$ws.Range("I3").Value = 0.006557053879060051
$ws.Range("J3").Value = 0.006557053879060051
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 27.685497
$ws.Range("N3").Value = 83.056491
$ws.Range("O3").Value = 0.420227262899125
$ws.Range("P3").Value = 0.4202272628991251
$ws.Range("Q3").Value = 3.172102732771
$ws.Range("R3").Value = 28.548924594939
$ws.Range("S3").Value = 0.002755452804279496
$ws.Range("T3").Value = 0.002755452804279496

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.1145763333333333
$ws.Range("H4").Value = 0.343729
$ws.Range("I4").Value = 0.006557053879060051
$ws.Range("J4").Value = 0.006557053879060051
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 4.125957666666666
$ws.Range("N4").Value = 12.377873
$ws.Range("O4").Value = 0.06262628758663766
$ws.Range("P4").Value = 0.06262628758663766
$ws.Range("Q4").Value = 0.4727371009352221
$ws.Range("R4").Value = 4.254633908417
$ws.Range("S4").Value = 0.0004106439419510927
$ws.Range("T4").Value = 0.0004106439419510927

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 9.390663666666667
$ws.Range("H5").Value = 28.171991
$ws.Range("I5").Value = 0.5374154140831726
$ws.Range("J5").Value = 0.5374154140831726
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 34.07074633333333
$ws.Range("N5").Value = 102.212239
$ws.Range("O5").Value = 0.5171464495142372
$ws.Range("P5").Value = 0.5171464495142373
$ws.Range("Q5").Value = 319.9469196886498
$ws.Range("R5").Value = 2879.522277197848
$ws.Range("S5").Value = 0.2779224733073363
$ws.Range("T5").Value = 0.2779224733073363

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 9.390663666666667
$ws.Range("H6").Value = 28.171991
$ws.Range("I6").Value = 0.5374154140831726
$ws.Range("J6").Value = 0.5374154140831726
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 27.685497
$ws.Range("N6").Value = 83.056491
$ws.Range("O6").Value = 0.420227262899125
$ws.Range("P6").Value = 0.4202272628991251
$ws.Range("Q6").Value = 259.985190771509
$ws.Range("R6").Value = 2339.866716943581
$ws.Range("S6").Value = 0.2258366084999715
$ws.Range("T6").Value = 0.2258366084999716

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 9.390663666666667
$ws.Range("H7").Value = 28.171991
$ws.Range("I7").Value = 0.5374154140831726
$ws.Range("J7").Value = 0.5374154140831726
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 4.125957666666666
$ws.Range("N7").Value = 12.377873
$ws.Range("O7").Value = 0.06262628758663766
$ws.Range("P7").Value = 0.06262628758663766
$ws.Range("Q7").Value = 38.74548075057145
$ws.Range("R7").Value = 348.709326755143
$ws.Range("S7").Value = 0.03365633227586473
$ws.Range("T7").Value = 0.03365633227586473

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 7.968512
$ws.Range("H8").Value = 23.905536
$ws.Range("I8").Value = 0.4560275320377672
$ws.Range("J8").Value = 0.4560275320377672
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 34.07074633333333
$ws.Range("N8").Value = 102.212239
$ws.Range("O8").Value = 0.5171464495142372
$ws.Range("P8").Value = 0.5171464495142373
$ws.Range("Q8").Value = 271.4931510061226
$ws.Range("R8").Value = 2443.438359055103
$ws.Range("S8").Value = 0.2358330190740714
$ws.Range("T8").Value = 0.2358330190740714

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 7.968512
$ws.Range("H9").Value = 23.905536
$ws.Range("I9").Value = 0.4560275320377672
$ws.Range("J9").Value = 0.4560275320377672
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 27.685497
$ws.Range("N9").Value = 83.056491
$ws.Range("O9").Value = 0.420227262899125
$ws.Range("P9").Value = 0.4202272628991251
$ws.Range("Q9").Value = 220.612215070464
$ws.Range("R9").Value = 1985.509935634176
$ws.Range("S9").Value = 0.191635201594874
$ws.Range("T9").Value = 0.191635201594874

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 7.968512
$ws.Range("H10").Value = 23.905536
$ws.Range("I10").Value = 0.4560275320377672
$ws.Range("J10").Value = 0.4560275320377672
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 4.125957666666666
$ws.Range("N10").Value = 12.377873
$ws.Range("O10").Value = 0.06262628758663766
$ws.Range("P10").Value = 0.06262628758663766
$ws.Range("Q10").Value = 32.87774317832533
$ws.Range("R10").Value = 295.899688604928
$ws.Range("S10").Value = 0.02855931136882183
$ws.Range("T10").Value = 0.02855931136882183
